$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "intent" column (B) for the newly added greeting/goodbye
# training rows (book_recomendation and main intent classification
# integration).

for ($r = 699; $r -le 712; $r++) {
    $ws.Cells.Item($r, 2).Value = "greeting"
}

for ($r = 713; $r -le 721; $r++) {
    $ws.Cells.Item($r, 2).Value = "goodbye"
}

# Update the view so it reflects where the user left off editing.
$excel.ActiveWindow.ScrollRow = 709
$ws.Range("B713:B721").Select()
